# Apply the commit's changes:
#  1. Rename the worksheet from "RGossF-HW40.xpc" to "RGossF"
#  2. Append a new data row (row 16) with the Gaussian-quadrature averaged
#     intensities for HKL index 14, using the same "HexGrid-60degTilt5degRes"
#     label as row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet/tab
$ws.Name = "RGossF"

# 2. Append the new row of data
$newRow = 16

# Copy the formatting (bold / border / centered) from the cell above (A15)
# onto the new A16 cell, then set its value.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 14

$ws.Cells.Item($newRow, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($newRow, 3).Value = 1.001031114751944
$ws.Cells.Item($newRow, 4).Value = 0.9872471667522974
$ws.Cells.Item($newRow, 5).Value = 1.000508245343851
$ws.Cells.Item($newRow, 6).Value = 1.001031114751944
$ws.Cells.Item($newRow, 7).Value = 0.9902760876083408
$ws.Cells.Item($newRow, 8).Value = 1.002854316891267
$ws.Cells.Item($newRow, 9).Value = 0.9994117647058823
$ws.Cells.Item($newRow, 10).Value = 0.9872471667522974
$ws.Cells.Item($newRow, 11).Value = 0.9938777060480739
$ws.Cells.Item($newRow, 12).Value = 0.9974544104000088
$ws.Cells.Item($newRow, 13).Value = 0.9968881160089302

$wb.Save()
